# Add a new worksheet "ODI Batting Extra" at the end of the workbook (after
# "ODI Bowling") and populate it with MATCH_CODE-keyed batting-extras data,
# matching the header styling already used on the other sheets.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# --- Header row ---
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Reuse the same bold / bordered / centered header style as the other sheets
# by copying the formatting from an existing header cell.
$srcWs = $wb.Worksheets.Item("ODI Bowling")
$srcWs.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Helper: write a value as TEXT (so digit-looking codes like "4402" are not
# auto-converted to numbers), without leaving a lasting number-format change
# on the cell.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# --- Data rows (keyed by MATCH_CODE, same order as the "ODI Batting" sheet) ---
Set-TextValue $ws.Cells.Item(2, 1) "4402"
Set-TextValue $ws.Cells.Item(2, 2) ""
Set-TextValue $ws.Cells.Item(2, 3) ""
Set-TextValue $ws.Cells.Item(2, 4) ""
Set-TextValue $ws.Cells.Item(2, 5) ""
$ws.Cells.Item(2, 6).Value = "NO"

Set-TextValue $ws.Cells.Item(3, 1) "4406"
Set-TextValue $ws.Cells.Item(3, 2) ""
Set-TextValue $ws.Cells.Item(3, 3) ""
Set-TextValue $ws.Cells.Item(3, 4) ""
Set-TextValue $ws.Cells.Item(3, 5) ""
$ws.Cells.Item(3, 6).Value = "NO"

Set-TextValue $ws.Cells.Item(4, 1) "4410"
$ws.Cells.Item(4, 2).Value = 2
Set-TextValue $ws.Cells.Item(4, 3) "0"
Set-TextValue $ws.Cells.Item(4, 4) "0"
Set-TextValue $ws.Cells.Item(4, 5) "0.34%"
$ws.Cells.Item(4, 6).Value = "NO"

Set-TextValue $ws.Cells.Item(5, 1) "4435"
Set-TextValue $ws.Cells.Item(5, 2) ""
Set-TextValue $ws.Cells.Item(5, 3) ""
Set-TextValue $ws.Cells.Item(5, 4) ""
Set-TextValue $ws.Cells.Item(5, 5) ""
$ws.Cells.Item(5, 6).Value = "NO"

Set-TextValue $ws.Cells.Item(6, 1) "4436"
Set-TextValue $ws.Cells.Item(6, 2) ""
Set-TextValue $ws.Cells.Item(6, 3) ""
Set-TextValue $ws.Cells.Item(6, 4) ""
Set-TextValue $ws.Cells.Item(6, 5) ""
$ws.Cells.Item(6, 6).Value = "NO"
